{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Rewrites the four bullet-style meeting notes into full sentences\n// describing the 30.11.2015 meeting, drops the trailing \"gesagt das...\"\n// paragraph, and adds a new blank paragraph before the closing line -\n// matching the commit \"Ferstigstellung der PM Arbeit, zudem alles in\n// PDF gewandelt\".\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraph 2 (0-based index): \"- Deinen text \u00fcberarbeitet und besprochen\"\nitems[2].insertText(\n  \"Dieses unplanm\u00e4\u00dfige Meeting am Montag, den 30.11.2015 wurde unter \" +\n  \"Beachtung eines Schwerpunktes durchgef\u00fchrt: die \u00dcberarbeitung der \" +\n  \"Analyse-Texte, die in der Alpha-Version vorliegen.\",\n  \"Replace\"\n);\n\n// Paragraph 3: \"- meinen text etwas \u00fcberarbeitet\"\nitems[3].insertText(\n  \"Der Analyse-Text, der Game Engines wurde korrigiert und verbessert. \" +\n  \"Dar\u00fcber hinaus wurde der Text in Microsoft Word verfasst. Da f\u00fcr die \" +\n  \"Dokumentation aber das Schreibprogramm Latex verwendet wird, mussten \" +\n  \"die Texte in Latex-Konvention angepasst werden.\",\n  \"Replace\"\n);\n\n// Paragraph 4: \"- Formatierung der texte\"\nitems[4].insertText(\n  \"Der Analyse-Text, ob Spiele etwas lehren, wurde bereits in Latex \" +\n  \"geschrieben. Aus diesem musste lediglich der Text korrigiert und \" +\n  \"verbessert werden.\",\n  \"Replace\"\n);\n\n// Paragraph 5: \"- \" (holds the _GoBack bookmark) -> becomes the closing\n// line; a new blank paragraph is inserted just before it.\nitems[5].insertParagraph(\"\", \"Before\");\nitems[5].insertText(\"Das n\u00e4chste Meeting findet am 01.12.2015 statt.\", \"Replace\");\n\n// Paragraph 6: \"gesagt das 01.12. n\u00e4chstes meeting\" -> removed entirely.\nitems[6].delete();\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Rewrites the four bullet-style meeting notes into full sentences\n# describing the 30.11.2015 meeting, drops the trailing \"gesagt das...\"\n# paragraph, and adds a new blank paragraph before the closing line -\n# matching the commit \"Ferstigstellung der PM Arbeit, zudem alles in\n# PDF gewandelt\".\n\n$d = $word.ActiveDocument\n\n# Paragraph 3 (1-based): \"- Deinen text \u00fcberarbeitet und besprochen\"\n# Re-grab the paragraph's own Range (start/end) and assign .Text on that\n# Range so the WHOLE paragraph (every run, including the proofErr-wrapped\n# \"text\" run) is replaced - assigning straight to Paragraphs(n).Range.Text\n# only clobbers the first run.\n$p3 = $d.Paragraphs(3)\n$d.Range($p3.Range.Start, $p3.Range.End).Text = \"Dieses unplanm\u00e4\u00dfige Meeting am Montag, den 30.11.2015 wurde unter Beachtung eines Schwerpunktes durchgef\u00fchrt: die \u00dcberarbeitung der Analyse-Texte, die in der Alpha-Version vorliegen.\"\n\n# Paragraph 4: \"- meinen text etwas \u00fcberarbeitet\"\n$p4 = $d.Paragraphs(4)\n$d.Range($p4.Range.Start, $p4.Range.End).Text = \"Der Analyse-Text, der Game Engines wurde korrigiert und verbessert. Dar\u00fcber hinaus wurde der Text in Microsoft Word verfasst. Da f\u00fcr die Dokumentation aber das Schreibprogramm Latex verwendet wird, mussten die Texte in Latex-Konvention angepasst werden.\"\n\n# Paragraph 5: \"- Formatierung der texte\"\n$p5 = $d.Paragraphs(5)\n$d.Range($p5.Range.Start, $p5.Range.End).Text = \"Der Analyse-Text, ob Spiele etwas lehren, wurde bereits in Latex geschrieben. Aus diesem musste lediglich der Text korrigiert und verbessert werden.\"\n\n# Paragraph 6: \"- \" (holds the _GoBack bookmark) -> becomes the closing\n# line; a new blank paragraph is inserted just before it, pushing the\n# bookmark paragraph down to index 7.\n$p6 = $d.Paragraphs(6)\n$p6.Range.InsertParagraphBefore()\n\n$p7 = $d.Paragraphs(7)\n$d.Range($p7.Range.Start, $p7.Range.End).Text = \"Das n\u00e4chste Meeting findet am 01.12.2015 statt.\"\n\n# Paragraph 8: \"gesagt das 01.12. n\u00e4chstes meeting\" -> removed entirely.\n$d.Paragraphs(8).Range.Delete()\n"}
